# Doing Updates for Financials
# Update the LOV (List Of Values) financial figures on the "LOV" sheet
# with the latest reported numbers (Income Statement, Balance Sheet and
# Cash Flow Statement sections).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Income Statement ---
$ws.Range("D8").Value = 96100
$ws.Range("E8").Value = 82500
$ws.Range("F8").Value = 67800
$ws.Range("D9").Value = 65900
$ws.Range("E9").Value = 57400
$ws.Range("F9").Value = 50100
$ws.Range("D10").Value = 30100
$ws.Range("E10").Value = 25000
$ws.Range("F10").Value = 17700
$ws.Range("D12").Value = 5800
$ws.Range("E12").Value = 3200
$ws.Range("F12").Value = 4300
$ws.Range("D15").Value = 3500
$ws.Range("E15").Value = 1400
$ws.Range("D17").Value = 101800
$ws.Range("E17").Value = 79500
$ws.Range("F17").Value = 66800
$ws.Range("D18").Value = -5700
$ws.Range("E18").Value = 3000
$ws.Range("F18").Value = 1000
$ws.Range("D21").Value = -2300
$ws.Range("E21").Value = 4300
$ws.Range("F21").Value = 1100
$ws.Range("E22").Value = 100
$ws.Range("D23").Value = -6300
$ws.Range("E23").Value = 2700
$ws.Range("E24").Value = 1200
$ws.Range("D26").Value = -6300
$ws.Range("E26").Value = 1500
$ws.Range("D27").Value = -6300
$ws.Range("E27").Value = 1500
$ws.Range("D33").Value = -6300
$ws.Range("D35").Value = -6300

# --- Balance Sheet ---
$ws.Range("D41").Value = 9200
$ws.Range("E41").Value = 9000
$ws.Range("F41").Value = 3400
$ws.Range("D42").Value = 3500
$ws.Range("E42").Value = 2800
$ws.Range("F42").Value = 2400
$ws.Range("D43").Value = 9000
$ws.Range("E43").Value = 5500
$ws.Range("F43").Value = 4600
$ws.Range("D45").Value = 3000
$ws.Range("E45").Value = 1900
$ws.Range("F45").Value = 2700
$ws.Range("D46").Value = 24700
$ws.Range("E46").Value = 19200
$ws.Range("F46").Value = 13100
$ws.Range("D48").Value = 2300
$ws.Range("E48").Value = 500
$ws.Range("D49").Value = 39400
$ws.Range("E49").Value = 10800
$ws.Range("D52").Value = 11100
$ws.Range("E52").Value = 11200
$ws.Range("F52").Value = 12200
$ws.Range("D54").Value = 77600
$ws.Range("E54").Value = 41800
$ws.Range("F54").Value = 25900
$ws.Range("D57").Value = 12900
$ws.Range("E57").Value = 6200
$ws.Range("F57").Value = 6100
$ws.Range("D58").Value = 6600
$ws.Range("D59").Value = 35500
$ws.Range("E59").Value = 26100
$ws.Range("F59").Value = 19700
$ws.Range("D60").Value = 54900
$ws.Range("E60").Value = 32300
$ws.Range("F60").Value = 25800
$ws.Range("E61").Value = 6600
$ws.Range("F61").Value = 29500
$ws.Range("E62").Value = 30600
$ws.Range("F62").Value = 200
$ws.Range("D66").Value = 55800
$ws.Range("E66").Value = 69600
$ws.Range("F66").Value = 55500
$ws.Range("D72").Value = 21400
$ws.Range("E72").Value = -27800
$ws.Range("F72").Value = -29700
$ws.Range("D76").Value = 21900
$ws.Range("E76").Value = -27700
$ws.Range("F76").Value = -29600

# --- Cash Flow Statement ---
$ws.Range("D81").Value = -6300
$ws.Range("D83").Value = 3500
$ws.Range("E83").Value = 1400
$ws.Range("D89").Value = -1300
$ws.Range("E89").Value = 7700
$ws.Range("F89").Value = 1600
$ws.Range("D91").Value = -2000
$ws.Range("E91").Value = -200
$ws.Range("D94").Value = 2700
$ws.Range("E94").Value = -8500
$ws.Range("E100").Value = 6600
$ws.Range("E102").Value = 5700
